# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update, Oct 21 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.599.24"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "1.600.11"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.46"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.87"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "1.829.29"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "1.597.78"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").Value = "29.614.57"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.73"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.91"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.99"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.25"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.16"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").Value = "1.426.84"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.02"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  +2.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.06"
$ws.Range("E41").Value = "  +7.00%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("E43").Value = "  +5.86%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.988"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.39"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "1.740.42"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.17"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  +3.96%  "
